# Auto-generated Excel COM-interop script
# Applies numeric updates to the "2024" (column K, and a couple column J) values
# across the Citywide Totals, By Neighborhood, and individual neighborhood sheets
# to add data for 2024-12-27, per the commit message "Add data for 2024-12-27".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7817
$ws.Range("K3").Value = 8100
$ws.Range("J4").Value = 1851
$ws.Range("K4").Value = 1705
$ws.Range("K6").Value = 9014
$ws.Range("J7").Value = 29321
$ws.Range("K7").Value = 27212

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 143
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 538
$ws.Range("K6").Value = 597
$ws.Range("K7").Value = 1778

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 202
$ws.Range("K6").Value = 137
$ws.Range("K7").Value = 579

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 286
$ws.Range("K3").Value = 405
$ws.Range("K4").Value = 57
$ws.Range("K6").Value = 362
$ws.Range("K7").Value = 1143

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 453

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 205
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 631

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 195
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 461

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 234
$ws.Range("K8").Value = 1778
$ws.Range("K10").Value = 165
$ws.Range("K18").Value = 185
$ws.Range("K19").Value = 781
$ws.Range("K20").Value = 671
$ws.Range("K22").Value = 84
$ws.Range("K24").Value = 88
$ws.Range("K27").Value = 261
$ws.Range("K29").Value = 1504
$ws.Range("K33").Value = 1143
$ws.Range("K36").Value = 354
$ws.Range("K41").Value = 179
$ws.Range("K42").Value = 1011
$ws.Range("K44").Value = 219
$ws.Range("K47").Value = 190
$ws.Range("K48").Value = 342
$ws.Range("K51").Value = 351
$ws.Range("K52").Value = 696
$ws.Range("K53").Value = 345
$ws.Range("K54").Value = 528
$ws.Range("J63").Value = 198
$ws.Range("K63").Value = 76
$ws.Range("K65").Value = 631
$ws.Range("K67").Value = 1063
$ws.Range("K73").Value = 241
$ws.Range("K75").Value = 89
$ws.Range("K77").Value = 179
$ws.Range("K83").Value = 579
$ws.Range("K85").Value = 1255
$ws.Range("K86").Value = 166
$ws.Range("K89").Value = 407
$ws.Range("K95").Value = 453
$ws.Range("K96").Value = 295
$ws.Range("K97").Value = 222
$ws.Range("K99").Value = 461
$ws.Range("J101").Value = 29321
$ws.Range("K101").Value = 27212

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 290
$ws.Range("K3").Value = 386
$ws.Range("K6").Value = 300
$ws.Range("K7").Value = 1063

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 284
$ws.Range("K7").Value = 528

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 424
$ws.Range("K7").Value = 1504

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 342

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 261
$ws.Range("K7").Value = 781

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 269
$ws.Range("K3").Value = 295
$ws.Range("K6").Value = 387
$ws.Range("K7").Value = 1011

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 295

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 209
$ws.Range("K6").Value = 196
$ws.Range("K7").Value = 671

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 185

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 132
$ws.Range("K7").Value = 354

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 71
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 126
$ws.Range("K7").Value = 407

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 63
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K4").Value = 40
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 414
$ws.Range("K3").Value = 436
$ws.Range("K6").Value = 308
$ws.Range("K7").Value = 1255

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 74
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 193
$ws.Range("K7").Value = 696
